$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying dataset (Covid-19 provincias Spain) was refreshed.
# Tenerife's total cases rose (1262 -> 1444), moving it above Asturias
# and Salamanca in the case-count sort order. Asturias and Salamanca keep
# their own figures but shift down one row (20->21, 21->22); Tenerife
# takes row 20 with its new figures.

# Row 20: now Tenerife
$ws.Range("A20").Value = "Tenerife"
$ws.Range("B20").Value = 1444
$ws.Range("C20").Value = 77
$ws.Range("D20").Value = 1241
$ws.Range("E20").Value = 68

# Row 21: now Asturias (unchanged figures, shifted down from row 20)
$ws.Range("A21").Value = "Asturias"
$ws.Range("B21").Value = 1384
$ws.Range("C21").Value = 135
$ws.Range("D21").Value = 1180
$ws.Range("E21").Value = 69

# Row 22: now Salamanca (unchanged figures, shifted down from row 21)
$ws.Range("A22").Value = "Salamanca"
$ws.Range("B22").Value = 1316
$ws.Range("C22").Value = 235
$ws.Range("D22").Value = 946
$ws.Range("E22").Value = 135

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 22:55"
